$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.305.52'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '2.684.41'
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'611.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').Value = "'160.44"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.80%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = "'0.595"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('D9').Value = "'0.127"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +10.54%  '
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('E13').Value = '  +24.91%  '
$ws.Range('D14').Value = "'30.56"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.78%  '
$ws.Range('D15').Value = '3.169.62'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').Value = '66.117.25'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '2.682.92'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = "'12.75"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').Value = "'4.91"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').Value = "'362.69"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').Value = "'7.50"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.41%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = "'70.21"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('D24').Value = "'9.78"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.45%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = "'0.0000108"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +17.67%  '
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').Value = "'1.69"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.76%  '
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('E28').Value = '  +5.81%  '
$ws.Range('D29').Value = "'8.21"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('E30').Value = '  +7.51%  '
$ws.Range('D31').Value = "'543.55"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').Value = "'6.62"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.59%  '
$ws.Range('D35').Value = "'5.59"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('D36').Value = "'0.436"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.96%  '
$ws.Range('D37').Value = "'20.84"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.42%  '
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').Value = "'163.12"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').Value = "'0.999"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = "'170.72"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = "'42.44"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').Value = "'4.28"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.24%  '
$ws.Range('D45').Value = "'2.37"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.97%  '
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').Value = "'23.27"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').Value = "'0.665"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('E49').Value = '  +5.58%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = "'0.0993"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'20.03"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.37%  '
